$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: refresh aggregate stats now that trade #7 has closed.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -0.23   # Total P&L %
$summary.Range("B6").Value = 7       # Total Trades
$summary.Range("B9").Value = 42.86   # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet: MarketMaking row picks up the new trade count /
# win rate too.
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 7        # Trades
$status.Range("G4").Value = 42.86    # Win Rate %

# ---------------------------------------------------------------------------
# Helper: append the newly closed trade (#7) as a new row on a trade log
# sheet, by duplicating the last existing row (so the date column keeps its
# original text formatting instead of being reinterpreted as a date serial)
# and then overwriting the cells that actually differ.
# ---------------------------------------------------------------------------
function Add-TradeSevenRow($ws) {
    $ws.Range("A7:Q7").Copy()
    $ws.Range("A8:Q8").PasteSpecial()

    $ws.Cells.Item(8, 1).Value = 7              # Trade #
    # B8 (Date) stays "2026-02-17", copied verbatim from row 7
    $ws.Cells.Item(8, 3).Value = "08:08:04"     # Time
    # D8 (Strategy) stays "MarketMaking", copied verbatim from row 7
    # E8 (Side) stays "DOWN", copied verbatim from row 7
    $ws.Cells.Item(8, 6).Value = 0.71           # Entry Price
    $ws.Cells.Item(8, 7).Value = 0.71           # Exit Price
    # H8 (Status) stays "CLOSED", copied verbatim from row 7
    $ws.Cells.Item(8, 9).Value = 0              # P&L %
    $ws.Cells.Item(8, 10).Value = 0             # P&L $
    # K8 (Capital After) stays 99.92, copied verbatim from row 7
    # L8 (Entry Slippage) stays 0, copied verbatim from row 7
    # M8 (Exit Slippage) stays 0, copied verbatim from row 7
    # N8 (Confidence) stays 0.6, copied verbatim from row 7
    # O8 (Entry Reason) stays "Normal spread capture: 19600 bps"
    # P8 (Exit Reason) stays "early_exit"
    # Q8 (Duration) stays 0.13, copied verbatim from row 7
}

# ---------------------------------------------------------------------------
# All Trades sheet: append trade #7.
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeSevenRow $allTrades

# ---------------------------------------------------------------------------
# MarketMaking sheet: append trade #7 (mirrors All Trades).
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeSevenRow $marketMaking
